$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rgValues = @(
    "389447-Amapá",
    "165808-Ssp/ap",
    "094524-Amapá",
    "080062-Amapá",
    "143983-Ap",
    "256017-Ssap",
    "7210826-Pará",
    "401687-Politec",
    "143876-Ap",
    "487237-Amapá",
    "529369-Politec",
    "332015-Ap",
    "481915-Amapá",
    "385779-Ap",
    "613689-Amapá",
    "187983-Amapá",
    "732139-Amapá",
    "646021-Ap",
    "504119-Ap",
    "331710-Ap",
    "376253-Ap",
    "469015-Amapá",
    "391493-Ssp",
    "164816-Amapá",
    "669892-Amapá",
    "586278-Ap",
    "673625-Amapá",
    "147286-Ap",
    "442179-Ap",
    "449442-Amapá",
    "625752-Siac",
    "470862-Amapá (Ap)",
    "610531-Amapá",
    "550792-Macapá",
    "505138-Amapá",
    "146434-Ap",
    "413039-Ap",
    "491458-Amapá",
    "409252-Ap",
    "1717-50-Amapa",
    "701003-Amapá",
    "483595-Amapá",
    "476827-Amapá",
    "352543-Amapá",
    "085302-Amapá",
    "667845537-São Paulo",
    "425079-Amapa",
    "293946-Amapá",
    "387936-Politec",
    "157822-Ap",
    "533558-Amapá",
    "553908-Ap",
    "173339-Ap",
    "580485-Amapá",
    "551400-Ap",
    "586482-Amapaense",
    "554951-Ap",
    "123873-Ap",
    "198142-Amapá",
    "528435-Ap",
    "7322152-Pará",
    "642338-Ap",
    "667380-Amapá",
    "135545-Ap",
    "594437-Amapá",
    "335471-Amapá",
    "524163-Amapá",
    "452028-Amapá",
    "460662-Ap",
    "554636-Siac",
    "296370-Amapá",
    "368815-Politec",
    "417036-Amapá",
    "6272243-Policia Civil",
    "647586-Ap",
    "195374-Ap",
    "431657-Ap",
    "549525-Ap",
    "857257-Politec",
    "900950-Solteira",
    "174930-Siac",
    "515252-Ap",
    "847010-Ptcdicc",
    "447448-Ap",
    "422461-Ap",
    "621546-Politec",
    "385407-Ap",
    "226048-Ap",
    "473664-Amapá",
    "577490-Amapá",
    "374203-Ap",
    "573416-Ap",
    "103905-Amapá",
    "555691-Amapá",
    "6243910-Pa",
    "522938-Amapá",
    "440227-Ap",
    "509370-Amapá",
    "593804-Amapá",
    "411401-Ap",
    "524014-Siac",
    "356207-Ap",
    "341370-Ap",
    "647940-Amapá",
    "164011-Ap",
    "520746-Ap",
    "156379-Amapá",
    "5081029-Segup/pa",
    "505162-Ap",
    "457881-Amapá",
    "358821-Amapá",
    "256016-Ssap",
    "559209-Politec",
    "672843-Amapá - Ap",
    "599345-Ap",
    "598818-Amapá",
    "904570-Politec",
    "594422-Amapá",
    "612823-Amapá",
    "467340-Ap",
    "688628-Amapá",
    "7107347-Pa",
    "650232-Amapá",
    "425087-Ap",
    "555952-Ap",
    "256472-Ap",
    "469693-Ap",
    "549526-Amapá",
    "337750-Ap",
    "510487-Ap",
    "366180-Amapá",
    "569357-Ap",
    "280699-Amapá",
    "131884-Amapá",
    "471864-Dptc",
    "564176-Amapá",
    "106179-Amapá",
    "15806860-5-Ssp/pr",
    "0738416620218-Maranhão",
    "340573-Ap",
    "328870-Amapá",
    "7902213-Policia Civil"
)

for ($i = 0; $i -lt $rgValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $rgValues[$i]
}

$ws.Cells.Item(1, 7).Value = "RG"

$ws.Columns.Item(6).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(7).EntireColumn.AutoFit() | Out-Null

$ws.Range("I5").Select() | Out-Null
